# Populate the "Port Configuration" sheet (sheet1) with the Meraki port
# table data that was added in this commit (rows 2-9), plus two trailing
# blank-but-styled rows (10-11), and move the active selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Port Configuration")

$HOSTNAME = "TEST-HOSTNAME"
$SERIAL   = "Q2GP-5GKZ-E8WY"

# One hashtable per data row (row numbers match the worksheet rows).
$rows = @(
  @{ Row=2; Port=1; Name="Staging"; Enabled=$true; Rstp=$true; StpGuard="disabled"; Poe=$false; Type="trunk"; Vlan=1; VoiceVlan=$null; Allowed="all" },
  @{ Row=3; Port=2; Name="Admin (v120), Voice (v300)"; Enabled=$true; Rstp=$true; StpGuard="disabled"; Poe=$true; Type="access"; Vlan=120; VoiceVlan=300; Allowed=$null },
  @{ Row=4; Port=3; Name="Admin (v120), Voice (v300)"; Enabled=$true; Rstp=$true; StpGuard="disabled"; Poe=$true; Type="access"; Vlan=120; VoiceVlan=300; Allowed=$null },
  @{ Row=5; Port=4; Name="Admin (v120), Voice (v300)"; Enabled=$true; Rstp=$true; StpGuard="disabled"; Poe=$true; Type="access"; Vlan=120; VoiceVlan=300; Allowed=$null },
  @{ Row=6; Port=5; Name="Trunk to MR-33"; Enabled=$true; Rstp=$false; StpGuard=$null; Poe=$true; Type="trunk"; Vlan=200; VoiceVlan=$null; Allowed="200,210,220,232,236,400" },
  @{ Row=7; Port=6; Name="Trunk to MR-33"; Enabled=$true; Rstp=$false; StpGuard=$null; Poe=$true; Type="trunk"; Vlan=200; VoiceVlan=$null; Allowed="200,210,220,232,236,400" },
  @{ Row=8; Port=7; Name="Trunk to MR-33"; Enabled=$true; Rstp=$false; StpGuard=$null; Poe=$true; Type="trunk"; Vlan=200; VoiceVlan=$null; Allowed="200,210,220,232,236,400" },
  @{ Row=9; Port=8; Name="SHUTDOWN - Not In Use"; Enabled=$false; Rstp=$null; StpGuard=$null; Poe=$null; Type=$null; Vlan=$null; VoiceVlan=$null; Allowed=$null }
)

foreach ($r in $rows) {
  $row = $r.Row

  # A: Hostname, B: Serial Number, C: Port Number, D: Name
  $ws.Cells.Item($row, 1).Value = $HOSTNAME
  $ws.Cells.Item($row, 2).Value = $SERIAL
  $ws.Cells.Item($row, 3).Value = $r.Port
  $ws.Cells.Item($row, 4).Value = $r.Name

  # F: Enabled, G: RSTP, H: STP Guard, I: PoE
  $ws.Cells.Item($row, 6).Value = $r.Enabled
  if ($null -ne $r.Rstp) { $ws.Cells.Item($row, 7).Value = $r.Rstp }
  if ($r.StpGuard) { $ws.Cells.Item($row, 8).Value = $r.StpGuard }
  if ($null -ne $r.Poe) { $ws.Cells.Item($row, 9).Value = $r.Poe }

  # J: Type, K: VLAN, L: Voice VLAN
  if ($r.Type) { $ws.Cells.Item($row, 10).Value = $r.Type }
  if ($r.Vlan) { $ws.Cells.Item($row, 11).Value = $r.Vlan }
  if ($r.VoiceVlan) { $ws.Cells.Item($row, 12).Value = $r.VoiceVlan }

  # M: Allowed VLANs for row 2 only entered here; the repeated "trunk to
  # MR-33" ports (rows 6-8) get theirs filled in afterwards (see below),
  # matching how the source sheet was authored (fill/paste down column M
  # as a separate pass once the rest of the table was complete).
  if ($row -eq 2 -and $r.Allowed) { $ws.Cells.Item($row, 13).Value = $r.Allowed }
}

foreach ($r in $rows) {
  if ($r.Row -ne 2 -and $r.Allowed) {
    $ws.Cells.Item($r.Row, 13).Value = $r.Allowed
  }
}

# Two trailing rows with only column A present, styled with an explicit
# black font color (new font/cellXf picked up automatically).
$ws.Range("A10:A11").Font.Color = 0

# Move the active selection, matching the post-edit worksheet view.
[void]$ws.Range("D5").Select()
